# Apply weekly update: shift rows 300-371 down by one (new_row[r] = old_row[r-1])
# and insert a brand-new observation into row 300, pushing the former
# last row (371) down into a newly appended row 372.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 372 is brand new (beyond the old A1:R371 used range) -> seed the
# row-constant columns that every data row shares.
$ws.Range("A372").Value = 4
$ws.Range("B372").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C372").Value = "Los Lagos"
$ws.Range("E372").Value = 10
$ws.Range("F372").Value = 100112021
$ws.Range("G372").Value = "Ají"
$ws.Range("R372").Value = "Hortaliza"

# Row 300
$ws.Range("D300").Value = 44964
$ws.Range("D300").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H300").Value = "Inferno"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 160
$ws.Range("K300").Value = 18000
$ws.Range("L300").Value = 20000
$ws.Range("M300").Value = 19000
$ws.Range("N300").Value = "`$/caja 10 kilos"
$ws.Range("O300").Value = "Región de Arica y Parinacota"
$ws.Range("P300").Value = 1900
$ws.Range("Q300").Value = 10

# Row 301
$ws.Range("D301").Value = 44900
$ws.Range("D301").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H301").Value = "Inferno"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 70
$ws.Range("K301").Value = 21000
$ws.Range("L301").Value = 21000
$ws.Range("M301").Value = 21000
$ws.Range("N301").Value = "`$/caja 10 kilos"
$ws.Range("O301").Value = "Región de Arica y Parinacota"
$ws.Range("P301").Value = 2100
$ws.Range("Q301").Value = 10

# Row 302
$ws.Range("D302").Value = 44476
$ws.Range("D302").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H302").Value = "Inferno"
$ws.Range("I302").Value = "Primera"
$ws.Range("J302").Value = 60
$ws.Range("K302").Value = 50000
$ws.Range("L302").Value = 50000
$ws.Range("M302").Value = 50000
$ws.Range("N302").Value = "`$/caja 12 kilos"
$ws.Range("O302").Value = "Región de Arica y Parinacota"
$ws.Range("P302").Value = 4167
$ws.Range("Q302").Value = 12

# Row 303
$ws.Range("D303").Value = 44579
$ws.Range("D303").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H303").Value = "Inferno"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 80
$ws.Range("K303").Value = 25000
$ws.Range("L303").Value = 25000
$ws.Range("M303").Value = 25000
$ws.Range("N303").Value = "`$/caja 15 kilos"
$ws.Range("O303").Value = "Región Metropolitana"
$ws.Range("P303").Value = 1667
$ws.Range("Q303").Value = 15

# Row 304
$ws.Range("D304").Value = 44711
$ws.Range("D304").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H304").Value = "Inferno"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 35
$ws.Range("K304").Value = 30000
$ws.Range("L304").Value = 30000
$ws.Range("M304").Value = 30000
$ws.Range("N304").Value = "`$/caja 12 kilos"
$ws.Range("O304").Value = "Región de Arica y Parinacota"
$ws.Range("P304").Value = 2500
$ws.Range("Q304").Value = 12

# Row 305
$ws.Range("D305").Value = 44161
$ws.Range("D305").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H305").Value = "Inferno"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 40
$ws.Range("K305").Value = 30000
$ws.Range("L305").Value = 30000
$ws.Range("M305").Value = 30000
$ws.Range("N305").Value = "`$/caja 12 kilos"
$ws.Range("O305").Value = "Región de Arica y Parinacota"
$ws.Range("P305").Value = 2500
$ws.Range("Q305").Value = 12

# Row 306
$ws.Range("D306").Value = 44882
$ws.Range("D306").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H306").Value = "Inferno"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 100
$ws.Range("K306").Value = 20000
$ws.Range("L306").Value = 21000
$ws.Range("M306").Value = 20500
$ws.Range("N306").Value = "`$/caja 10 kilos"
$ws.Range("O306").Value = "Región de Arica y Parinacota"
$ws.Range("P306").Value = 2050
$ws.Range("Q306").Value = 10

# Row 307
$ws.Range("D307").Value = 44264
$ws.Range("D307").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H307").Value = "Chilena(o)"
$ws.Range("I307").Value = "Primera"
$ws.Range("J307").Value = 40
$ws.Range("K307").Value = 20000
$ws.Range("L307").Value = 20000
$ws.Range("M307").Value = 20000
$ws.Range("N307").Value = "`$/saco 25 kilos"
$ws.Range("O307").Value = "Región Metropolitana"
$ws.Range("P307").Value = 800
$ws.Range("Q307").Value = 25

# Row 308
$ws.Range("D308").Value = 44264
$ws.Range("D308").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H308").Value = "Inferno"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 80
$ws.Range("K308").Value = 18000
$ws.Range("L308").Value = 18000
$ws.Range("M308").Value = 18000
$ws.Range("N308").Value = "`$/caja 14 kilos"
$ws.Range("O308").Value = "Limache"
$ws.Range("P308").Value = 1286
$ws.Range("Q308").Value = 14

# Row 309
$ws.Range("D309").Value = 44407
$ws.Range("D309").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H309").Value = "Inferno"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 80
$ws.Range("K309").Value = 32000
$ws.Range("L309").Value = 32000
$ws.Range("M309").Value = 32000
$ws.Range("N309").Value = "`$/caja 12 kilos"
$ws.Range("O309").Value = "Región de Arica y Parinacota"
$ws.Range("P309").Value = 2667
$ws.Range("Q309").Value = 12

# Row 310
$ws.Range("D310").Value = 44407
$ws.Range("D310").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H310").Value = "Inferno"
$ws.Range("I310").Value = "Segunda"
$ws.Range("J310").Value = 80
$ws.Range("K310").Value = 26000
$ws.Range("L310").Value = 26000
$ws.Range("M310").Value = 26000
$ws.Range("N310").Value = "`$/caja 12 kilos"
$ws.Range("O310").Value = "Región de Arica y Parinacota"
$ws.Range("P310").Value = 2167
$ws.Range("Q310").Value = 12

# Row 311
$ws.Range("D311").Value = 44250
$ws.Range("D311").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H311").Value = "Inferno"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 150
$ws.Range("K311").Value = 15000
$ws.Range("L311").Value = 15000
$ws.Range("M311").Value = 15000
$ws.Range("N311").Value = "`$/caja 14 kilos"
$ws.Range("O311").Value = "Provincia de Quillota"
$ws.Range("P311").Value = 1071
$ws.Range("Q311").Value = 14

# Row 312
$ws.Range("D312").Value = 44215
$ws.Range("D312").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H312").Value = "Inferno"
$ws.Range("I312").Value = "Primera"
$ws.Range("J312").Value = 120
$ws.Range("K312").Value = 35000
$ws.Range("L312").Value = 35000
$ws.Range("M312").Value = 35000
$ws.Range("N312").Value = "`$/caja 14 kilos"
$ws.Range("O312").Value = "Provincia de Quillota"
$ws.Range("P312").Value = 2500
$ws.Range("Q312").Value = 14

# Row 313
$ws.Range("D313").Value = 44782
$ws.Range("D313").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H313").Value = "Inferno"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 140
$ws.Range("K313").Value = 21000
$ws.Range("L313").Value = 21000
$ws.Range("M313").Value = 21000
$ws.Range("N313").Value = "`$/caja 12 kilos"
$ws.Range("O313").Value = "Región de Arica y Parinacota"
$ws.Range("P313").Value = 1750
$ws.Range("Q313").Value = 12

# Row 314
$ws.Range("D314").Value = 44754
$ws.Range("D314").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H314").Value = "Inferno"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 150
$ws.Range("K314").Value = 21000
$ws.Range("L314").Value = 21000
$ws.Range("M314").Value = 21000
$ws.Range("N314").Value = "`$/caja 12 kilos"
$ws.Range("O314").Value = "Región de Arica y Parinacota"
$ws.Range("P314").Value = 1750
$ws.Range("Q314").Value = 12

# Row 315
$ws.Range("D315").Value = 44694
$ws.Range("D315").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H315").Value = "Inferno"
$ws.Range("I315").Value = "Primera"
$ws.Range("J315").Value = 80
$ws.Range("K315").Value = 31000
$ws.Range("L315").Value = 31000
$ws.Range("M315").Value = 31000
$ws.Range("N315").Value = "`$/caja 12 kilos"
$ws.Range("O315").Value = "Región de Arica y Parinacota"
$ws.Range("P315").Value = 2583
$ws.Range("Q315").Value = 12

# Row 316
$ws.Range("D316").Value = 44694
$ws.Range("D316").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H316").Value = "Inferno"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 70
$ws.Range("K316").Value = 34000
$ws.Range("L316").Value = 34000
$ws.Range("M316").Value = 34000
$ws.Range("N316").Value = "`$/caja 15 kilos"
$ws.Range("O316").Value = "Provincia de Huasco"
$ws.Range("P316").Value = 2267
$ws.Range("Q316").Value = 15

# Row 317
$ws.Range("D317").Value = 44663
$ws.Range("D317").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H317").Value = "Cristal"
$ws.Range("I317").Value = "Primera"
$ws.Range("J317").Value = 60
$ws.Range("K317").Value = 20000
$ws.Range("L317").Value = 20000
$ws.Range("M317").Value = 20000
$ws.Range("N317").Value = "`$/saco 25 kilos"
$ws.Range("O317").Value = "Región del Maule"
$ws.Range("P317").Value = 800
$ws.Range("Q317").Value = 25

# Row 318
$ws.Range("D318").Value = 44663
$ws.Range("D318").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H318").Value = "Inferno"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 120
$ws.Range("K318").Value = 30000
$ws.Range("L318").Value = 30000
$ws.Range("M318").Value = 30000
$ws.Range("N318").Value = "`$/caja 15 kilos"
$ws.Range("O318").Value = "Provincia de Quillota"
$ws.Range("P318").Value = 2000
$ws.Range("Q318").Value = 15

# Row 319
$ws.Range("D319").Value = 44399
$ws.Range("D319").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H319").Value = "Inferno"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 70
$ws.Range("K319").Value = 32000
$ws.Range("L319").Value = 32000
$ws.Range("M319").Value = 32000
$ws.Range("N319").Value = "`$/caja 12 kilos"
$ws.Range("O319").Value = "Región de Arica y Parinacota"
$ws.Range("P319").Value = 2667
$ws.Range("Q319").Value = 12

# Row 320
$ws.Range("D320").Value = 44266
$ws.Range("D320").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H320").Value = "Chilena(o)"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 20
$ws.Range("K320").Value = 20000
$ws.Range("L320").Value = 20000
$ws.Range("M320").Value = 20000
$ws.Range("N320").Value = "`$/saco 25 kilos"
$ws.Range("O320").Value = "Región Metropolitana"
$ws.Range("P320").Value = 800
$ws.Range("Q320").Value = 25

# Row 321
$ws.Range("D321").Value = 44446
$ws.Range("D321").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H321").Value = "Inferno"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 70
$ws.Range("K321").Value = 45000
$ws.Range("L321").Value = 45000
$ws.Range("M321").Value = 45000
$ws.Range("N321").Value = "`$/caja 12 kilos"
$ws.Range("O321").Value = "Región de Arica y Parinacota"
$ws.Range("P321").Value = 3750
$ws.Range("Q321").Value = 12

# Row 322
$ws.Range("D322").Value = 44446
$ws.Range("D322").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H322").Value = "Inferno"
$ws.Range("I322").Value = "Segunda"
$ws.Range("J322").Value = 70
$ws.Range("K322").Value = 35000
$ws.Range("L322").Value = 35000
$ws.Range("M322").Value = 35000
$ws.Range("N322").Value = "`$/caja 12 kilos"
$ws.Range("O322").Value = "Región de Arica y Parinacota"
$ws.Range("P322").Value = 2917
$ws.Range("Q322").Value = 12

# Row 323
$ws.Range("D323").Value = 44540
$ws.Range("D323").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H323").Value = "Inferno"
$ws.Range("I323").Value = "Primera"
$ws.Range("J323").Value = 160
$ws.Range("K323").Value = 18000
$ws.Range("L323").Value = 18000
$ws.Range("M323").Value = 18000
$ws.Range("N323").Value = "`$/caja 12 kilos"
$ws.Range("O323").Value = "Región de Arica y Parinacota"
$ws.Range("P323").Value = 1500
$ws.Range("Q323").Value = 12

# Row 324
$ws.Range("D324").Value = 44559
$ws.Range("D324").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H324").Value = "Inferno"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 140
$ws.Range("K324").Value = 24500
$ws.Range("L324").Value = 25000
$ws.Range("M324").Value = 24750
$ws.Range("N324").Value = "`$/caja 12 kilos"
$ws.Range("O324").Value = "Región de Arica y Parinacota"
$ws.Range("P324").Value = 2062
$ws.Range("Q324").Value = 12

# Row 325
$ws.Range("D325").Value = 44690
$ws.Range("D325").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H325").Value = "Inferno"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 20
$ws.Range("K325").Value = 31000
$ws.Range("L325").Value = 31000
$ws.Range("M325").Value = 31000
$ws.Range("N325").Value = "`$/caja 12 kilos"
$ws.Range("O325").Value = "Región de Arica y Parinacota"
$ws.Range("P325").Value = 2583
$ws.Range("Q325").Value = 12

# Row 326
$ws.Range("D326").Value = 44904
$ws.Range("D326").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H326").Value = "Inferno"
$ws.Range("I326").Value = "Primera"
$ws.Range("J326").Value = 240
$ws.Range("K326").Value = 19000
$ws.Range("L326").Value = 20000
$ws.Range("M326").Value = 19500
$ws.Range("N326").Value = "`$/caja 10 kilos"
$ws.Range("O326").Value = "Región de Arica y Parinacota"
$ws.Range("P326").Value = 1950
$ws.Range("Q326").Value = 10

# Row 327
$ws.Range("D327").Value = 44897
$ws.Range("D327").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H327").Value = "Inferno"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 100
$ws.Range("K327").Value = 21000
$ws.Range("L327").Value = 21000
$ws.Range("M327").Value = 21000
$ws.Range("N327").Value = "`$/caja 10 kilos"
$ws.Range("O327").Value = "Región de Arica y Parinacota"
$ws.Range("P327").Value = 2100
$ws.Range("Q327").Value = 10

# Row 328
$ws.Range("D328").Value = 44607
$ws.Range("D328").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H328").Value = "Inferno"
$ws.Range("I328").Value = "Primera"
$ws.Range("J328").Value = 150
$ws.Range("K328").Value = 18000
$ws.Range("L328").Value = 18000
$ws.Range("M328").Value = 18000
$ws.Range("N328").Value = "`$/caja 12 kilos"
$ws.Range("O328").Value = "Región de Arica y Parinacota"
$ws.Range("P328").Value = 1500
$ws.Range("Q328").Value = 12

# Row 329
$ws.Range("D329").Value = 44529
$ws.Range("D329").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H329").Value = "Inferno"
$ws.Range("I329").Value = "Primera"
$ws.Range("J329").Value = 40
$ws.Range("K329").Value = 22000
$ws.Range("L329").Value = 22000
$ws.Range("M329").Value = 22000
$ws.Range("N329").Value = "`$/caja 12 kilos"
$ws.Range("O329").Value = "Región de Arica y Parinacota"
$ws.Range("P329").Value = 1833
$ws.Range("Q329").Value = 12

# Row 330
$ws.Range("D330").Value = 44911
$ws.Range("D330").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H330").Value = "Inferno"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 180
$ws.Range("K330").Value = 19000
$ws.Range("L330").Value = 19000
$ws.Range("M330").Value = 19000
$ws.Range("N330").Value = "`$/caja 10 kilos"
$ws.Range("O330").Value = "Región de Arica y Parinacota"
$ws.Range("P330").Value = 1900
$ws.Range("Q330").Value = 10

# Row 331
$ws.Range("D331").Value = 44901
$ws.Range("D331").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H331").Value = "Inferno"
$ws.Range("I331").Value = "Primera"
$ws.Range("J331").Value = 180
$ws.Range("K331").Value = 20000
$ws.Range("L331").Value = 20000
$ws.Range("M331").Value = 20000
$ws.Range("N331").Value = "`$/caja 10 kilos"
$ws.Range("O331").Value = "Región de Arica y Parinacota"
$ws.Range("P331").Value = 2000
$ws.Range("Q331").Value = 10

# Row 332
$ws.Range("D332").Value = 44946
$ws.Range("D332").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H332").Value = "Inferno"
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 200
$ws.Range("K332").Value = 18000
$ws.Range("L332").Value = 20000
$ws.Range("M332").Value = 19000
$ws.Range("N332").Value = "`$/caja 10 kilos"
$ws.Range("O332").Value = "Región de Arica y Parinacota"
$ws.Range("P332").Value = 1900
$ws.Range("Q332").Value = 10

# Row 333
$ws.Range("D333").Value = 44467
$ws.Range("D333").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H333").Value = "Inferno"
$ws.Range("I333").Value = "Primera"
$ws.Range("J333").Value = 120
$ws.Range("K333").Value = 48000
$ws.Range("L333").Value = 48000
$ws.Range("M333").Value = 48000
$ws.Range("N333").Value = "`$/caja 12 kilos"
$ws.Range("O333").Value = "Región de Arica y Parinacota"
$ws.Range("P333").Value = 4000
$ws.Range("Q333").Value = 12

# Row 334
$ws.Range("D334").Value = 44628
$ws.Range("D334").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H334").Value = "Inferno"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 120
$ws.Range("K334").Value = 20000
$ws.Range("L334").Value = 20000
$ws.Range("M334").Value = 20000
$ws.Range("N334").Value = "`$/caja 15 kilos"
$ws.Range("O334").Value = "Región Metropolitana"
$ws.Range("P334").Value = 1333
$ws.Range("Q334").Value = 15

# Row 335
$ws.Range("D335").Value = 44644
$ws.Range("D335").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H335").Value = "Cristal"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 20
$ws.Range("K335").Value = 20000
$ws.Range("L335").Value = 20000
$ws.Range("M335").Value = 20000
$ws.Range("N335").Value = "`$/saco 25 kilos"
$ws.Range("O335").Value = "Región del Maule"
$ws.Range("P335").Value = 800
$ws.Range("Q335").Value = 25

# Row 336
$ws.Range("D336").Value = 44224
$ws.Range("D336").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H336").Value = "Inferno"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 60
$ws.Range("K336").Value = 30000
$ws.Range("L336").Value = 30000
$ws.Range("M336").Value = 30000
$ws.Range("N336").Value = "`$/caja 14 kilos"
$ws.Range("O336").Value = "Provincia de Quillota"
$ws.Range("P336").Value = 2143
$ws.Range("Q336").Value = 14

# Row 337
$ws.Range("D337").Value = 44505
$ws.Range("D337").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H337").Value = "Inferno"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 70
$ws.Range("K337").Value = 33000
$ws.Range("L337").Value = 33000
$ws.Range("M337").Value = 33000
$ws.Range("N337").Value = "`$/caja 12 kilos"
$ws.Range("O337").Value = "Región de Arica y Parinacota"
$ws.Range("P337").Value = 2750
$ws.Range("Q337").Value = 12

# Row 338
$ws.Range("D338").Value = 44505
$ws.Range("D338").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H338").Value = "Inferno"
$ws.Range("I338").Value = "Segunda"
$ws.Range("J338").Value = 70
$ws.Range("K338").Value = 28000
$ws.Range("L338").Value = 28000
$ws.Range("M338").Value = 28000
$ws.Range("N338").Value = "`$/caja 12 kilos"
$ws.Range("O338").Value = "Región de Arica y Parinacota"
$ws.Range("P338").Value = 2333
$ws.Range("Q338").Value = 12

# Row 339
$ws.Range("D339").Value = 44637
$ws.Range("D339").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H339").Value = "Inferno"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 60
$ws.Range("K339").Value = 23000
$ws.Range("L339").Value = 23000
$ws.Range("M339").Value = 23000
$ws.Range("N339").Value = "`$/caja 15 kilos"
$ws.Range("O339").Value = "Provincia de Quillota"
$ws.Range("P339").Value = 1533
$ws.Range("Q339").Value = 15

# Row 340
$ws.Range("D340").Value = 44239
$ws.Range("D340").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H340").Value = "Inferno"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 150
$ws.Range("K340").Value = 15000
$ws.Range("L340").Value = 15000
$ws.Range("M340").Value = 15000
$ws.Range("N340").Value = "`$/caja 12 kilos"
$ws.Range("O340").Value = "Región de Arica y Parinacota"
$ws.Range("P340").Value = 1250
$ws.Range("Q340").Value = 12

# Row 341
$ws.Range("D341").Value = 44855
$ws.Range("D341").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H341").Value = "Inferno"
$ws.Range("I341").Value = "Primera"
$ws.Range("J341").Value = 90
$ws.Range("K341").Value = 25000
$ws.Range("L341").Value = 25000
$ws.Range("M341").Value = 25000
$ws.Range("N341").Value = "`$/caja 10 kilos"
$ws.Range("O341").Value = "Región de Arica y Parinacota"
$ws.Range("P341").Value = 2500
$ws.Range("Q341").Value = 10

# Row 342
$ws.Range("D342").Value = 44855
$ws.Range("D342").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H342").Value = "Inferno"
$ws.Range("I342").Value = "Segunda"
$ws.Range("J342").Value = 90
$ws.Range("K342").Value = 20000
$ws.Range("L342").Value = 20000
$ws.Range("M342").Value = 20000
$ws.Range("N342").Value = "`$/caja 10 kilos"
$ws.Range("O342").Value = "Región de Arica y Parinacota"
$ws.Range("P342").Value = 2000
$ws.Range("Q342").Value = 10

# Row 343
$ws.Range("D343").Value = 44616
$ws.Range("D343").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H343").Value = "Inferno"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 60
$ws.Range("K343").Value = 20000
$ws.Range("L343").Value = 20000
$ws.Range("M343").Value = 20000
$ws.Range("N343").Value = "`$/caja 15 kilos"
$ws.Range("O343").Value = "Región Metropolitana"
$ws.Range("P343").Value = 1333
$ws.Range("Q343").Value = 15

# Row 344
$ws.Range("D344").Value = 44581
$ws.Range("D344").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H344").Value = "Inferno"
$ws.Range("I344").Value = "Primera"
$ws.Range("J344").Value = 70
$ws.Range("K344").Value = 18000
$ws.Range("L344").Value = 18000
$ws.Range("M344").Value = 18000
$ws.Range("N344").Value = "`$/caja 12 kilos"
$ws.Range("O344").Value = "Región de Arica y Parinacota"
$ws.Range("P344").Value = 1500
$ws.Range("Q344").Value = 12

# Row 345
$ws.Range("D345").Value = 44271
$ws.Range("D345").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H345").Value = "Inferno"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 110
$ws.Range("K345").Value = 18000
$ws.Range("L345").Value = 18000
$ws.Range("M345").Value = 18000
$ws.Range("N345").Value = "`$/caja 14 kilos"
$ws.Range("O345").Value = "Provincia de Quillota"
$ws.Range("P345").Value = 1286
$ws.Range("Q345").Value = 14

# Row 346
$ws.Range("D346").Value = 44908
$ws.Range("D346").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H346").Value = "Inferno"
$ws.Range("I346").Value = "Primera"
$ws.Range("J346").Value = 180
$ws.Range("K346").Value = 19000
$ws.Range("L346").Value = 19000
$ws.Range("M346").Value = 19000
$ws.Range("N346").Value = "`$/caja 10 kilos"
$ws.Range("O346").Value = "Región de Arica y Parinacota"
$ws.Range("P346").Value = 1900
$ws.Range("Q346").Value = 10

# Row 347
$ws.Range("D347").Value = 44259
$ws.Range("D347").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H347").Value = "Inferno"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 30
$ws.Range("K347").Value = 18000
$ws.Range("L347").Value = 18000
$ws.Range("M347").Value = 18000
$ws.Range("N347").Value = "`$/caja 14 kilos"
$ws.Range("O347").Value = "Provincia de Quillota"
$ws.Range("P347").Value = 1286
$ws.Range("Q347").Value = 14

# Row 348
$ws.Range("D348").Value = 44252
$ws.Range("D348").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H348").Value = "Inferno"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 40
$ws.Range("K348").Value = 16000
$ws.Range("L348").Value = 16000
$ws.Range("M348").Value = 16000
$ws.Range("N348").Value = "`$/caja 14 kilos"
$ws.Range("O348").Value = "Provincia de Quillota"
$ws.Range("P348").Value = 1143
$ws.Range("Q348").Value = 14

# Row 349
$ws.Range("D349").Value = 44243
$ws.Range("D349").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H349").Value = "Chilena(o)"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 40
$ws.Range("K349").Value = 18000
$ws.Range("L349").Value = 18000
$ws.Range("M349").Value = 18000
$ws.Range("N349").Value = "`$/saco 25 kilos"
$ws.Range("O349").Value = "Región Metropolitana"
$ws.Range("P349").Value = 720
$ws.Range("Q349").Value = 25

# Row 350
$ws.Range("D350").Value = 44243
$ws.Range("D350").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H350").Value = "Inferno"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 120
$ws.Range("K350").Value = 15000
$ws.Range("L350").Value = 15000
$ws.Range("M350").Value = 15000
$ws.Range("N350").Value = "`$/caja 14 kilos"
$ws.Range("O350").Value = "Provincia de Quillota"
$ws.Range("P350").Value = 1071
$ws.Range("Q350").Value = 14

# Row 351
$ws.Range("D351").Value = 44539
$ws.Range("D351").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H351").Value = "Inferno"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 70
$ws.Range("K351").Value = 18000
$ws.Range("L351").Value = 20000
$ws.Range("M351").Value = 19143
$ws.Range("N351").Value = "`$/caja 12 kilos"
$ws.Range("O351").Value = "Región de Arica y Parinacota"
$ws.Range("P351").Value = 1595
$ws.Range("Q351").Value = 12

# Row 352
$ws.Range("D352").Value = 44826
$ws.Range("D352").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H352").Value = "Inferno"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 90
$ws.Range("K352").Value = 23000
$ws.Range("L352").Value = 23000
$ws.Range("M352").Value = 23000
$ws.Range("N352").Value = "`$/caja 10 kilos"
$ws.Range("O352").Value = "Región de Arica y Parinacota"
$ws.Range("P352").Value = 2300
$ws.Range("Q352").Value = 10

# Row 353
$ws.Range("D353").Value = 44757
$ws.Range("D353").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H353").Value = "Inferno"
$ws.Range("I353").Value = "Primera"
$ws.Range("J353").Value = 150
$ws.Range("K353").Value = 21000
$ws.Range("L353").Value = 21000
$ws.Range("M353").Value = 21000
$ws.Range("N353").Value = "`$/caja 12 kilos"
$ws.Range("O353").Value = "Región de Arica y Parinacota"
$ws.Range("P353").Value = 1750
$ws.Range("Q353").Value = 12

# Row 354
$ws.Range("D354").Value = 44838
$ws.Range("D354").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H354").Value = "Americana (o)"
$ws.Range("I354").Value = "Primera"
$ws.Range("J354").Value = 90
$ws.Range("K354").Value = 26000
$ws.Range("L354").Value = 26000
$ws.Range("M354").Value = 26000
$ws.Range("N354").Value = "`$/caja 10 kilos"
$ws.Range("O354").Value = "Región de Arica y Parinacota"
$ws.Range("P354").Value = 2600
$ws.Range("Q354").Value = 10

# Row 355
$ws.Range("D355").Value = 44838
$ws.Range("D355").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H355").Value = "Americana (o)"
$ws.Range("I355").Value = "Segunda"
$ws.Range("J355").Value = 90
$ws.Range("K355").Value = 22000
$ws.Range("L355").Value = 22000
$ws.Range("M355").Value = 22000
$ws.Range("N355").Value = "`$/caja 10 kilos"
$ws.Range("O355").Value = "Región de Arica y Parinacota"
$ws.Range("P355").Value = 2200
$ws.Range("Q355").Value = 10

# Row 356
$ws.Range("D356").Value = 44229
$ws.Range("D356").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H356").Value = "Inferno"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 90
$ws.Range("K356").Value = 22000
$ws.Range("L356").Value = 22000
$ws.Range("M356").Value = 22000
$ws.Range("N356").Value = "`$/caja 12 kilos"
$ws.Range("O356").Value = "Región de Arica y Parinacota"
$ws.Range("P356").Value = 1833
$ws.Range("Q356").Value = 12

# Row 357
$ws.Range("D357").Value = 44229
$ws.Range("D357").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H357").Value = "Inferno"
$ws.Range("I357").Value = "Primera"
$ws.Range("J357").Value = 80
$ws.Range("K357").Value = 30000
$ws.Range("L357").Value = 30000
$ws.Range("M357").Value = 30000
$ws.Range("N357").Value = "`$/caja 14 kilos"
$ws.Range("O357").Value = "Provincia de Quillota"
$ws.Range("P357").Value = 2143
$ws.Range("Q357").Value = 14

# Row 358
$ws.Range("D358").Value = 44320
$ws.Range("D358").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H358").Value = "Chilena(o)"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 70
$ws.Range("K358").Value = 30000
$ws.Range("L358").Value = 30000
$ws.Range("M358").Value = 30000
$ws.Range("N358").Value = "`$/saco 25 kilos"
$ws.Range("O358").Value = "Región del Maule"
$ws.Range("P358").Value = 1200
$ws.Range("Q358").Value = 25

# Row 359
$ws.Range("D359").Value = 44320
$ws.Range("D359").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H359").Value = "Inferno"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 60
$ws.Range("K359").Value = 17000
$ws.Range("L359").Value = 17000
$ws.Range("M359").Value = 17000
$ws.Range("N359").Value = "`$/caja 15 kilos"
$ws.Range("O359").Value = "Provincia de Limarí"
$ws.Range("P359").Value = 1133
$ws.Range("Q359").Value = 15

# Row 360
$ws.Range("D360").Value = 44320
$ws.Range("D360").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H360").Value = "Inferno"
$ws.Range("I360").Value = "Segunda"
$ws.Range("J360").Value = 60
$ws.Range("K360").Value = 15000
$ws.Range("L360").Value = 15000
$ws.Range("M360").Value = 15000
$ws.Range("N360").Value = "`$/caja 15 kilos"
$ws.Range("O360").Value = "Provincia de Limarí"
$ws.Range("P360").Value = 1000
$ws.Range("Q360").Value = 15

# Row 361
$ws.Range("D361").Value = 44371
$ws.Range("D361").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H361").Value = "Inferno"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 60
$ws.Range("K361").Value = 31000
$ws.Range("L361").Value = 32000
$ws.Range("M361").Value = 31500
$ws.Range("N361").Value = "`$/caja 12 kilos"
$ws.Range("O361").Value = "Región de Arica y Parinacota"
$ws.Range("P361").Value = 2625
$ws.Range("Q361").Value = 12

# Row 362
$ws.Range("D362").Value = 44302
$ws.Range("D362").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H362").Value = "Inferno"
$ws.Range("I362").Value = "Segunda"
$ws.Range("J362").Value = 90
$ws.Range("K362").Value = 17000
$ws.Range("L362").Value = 17000
$ws.Range("M362").Value = 17000
$ws.Range("N362").Value = "`$/caja 15 kilos"
$ws.Range("O362").Value = "Región Metropolitana"
$ws.Range("P362").Value = 1133
$ws.Range("Q362").Value = 15

# Row 363
$ws.Range("D363").Value = 44592
$ws.Range("D363").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H363").Value = "Inferno"
$ws.Range("I363").Value = "Primera"
$ws.Range("J363").Value = 35
$ws.Range("K363").Value = 17000
$ws.Range("L363").Value = 17000
$ws.Range("M363").Value = 17000
$ws.Range("N363").Value = "`$/caja 12 kilos"
$ws.Range("O363").Value = "Región de Arica y Parinacota"
$ws.Range("P363").Value = 1417
$ws.Range("Q363").Value = 12

# Row 364
$ws.Range("D364").Value = 44225
$ws.Range("D364").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H364").Value = "Inferno"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 100
$ws.Range("K364").Value = 25000
$ws.Range("L364").Value = 30000
$ws.Range("M364").Value = 28000
$ws.Range("N364").Value = "`$/caja 14 kilos"
$ws.Range("O364").Value = "Provincia de Quillota"
$ws.Range("P364").Value = 2000
$ws.Range("Q364").Value = 14

# Row 365
$ws.Range("D365").Value = 44810
$ws.Range("D365").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H365").Value = "Inferno"
$ws.Range("I365").Value = "Primera"
$ws.Range("J365").Value = 100
$ws.Range("K365").Value = 25000
$ws.Range("L365").Value = 25000
$ws.Range("M365").Value = 25000
$ws.Range("N365").Value = "`$/caja 10 kilos"
$ws.Range("O365").Value = "Región de Arica y Parinacota"
$ws.Range("P365").Value = 2500
$ws.Range("Q365").Value = 10

# Row 366
$ws.Range("D366").Value = 44810
$ws.Range("D366").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H366").Value = "Inferno"
$ws.Range("I366").Value = "Segunda"
$ws.Range("J366").Value = 100
$ws.Range("K366").Value = 20000
$ws.Range("L366").Value = 20000
$ws.Range("M366").Value = 20000
$ws.Range("N366").Value = "`$/caja 10 kilos"
$ws.Range("O366").Value = "Región de Arica y Parinacota"
$ws.Range("P366").Value = 2000
$ws.Range("Q366").Value = 10

# Row 367
$ws.Range("D367").Value = 44175
$ws.Range("D367").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H367").Value = "Inferno"
$ws.Range("I367").Value = "Primera"
$ws.Range("J367").Value = 60
$ws.Range("K367").Value = 29000
$ws.Range("L367").Value = 30000
$ws.Range("M367").Value = 29500
$ws.Range("N367").Value = "`$/caja 12 kilos"
$ws.Range("O367").Value = "Región de Arica y Parinacota"
$ws.Range("P367").Value = 2458
$ws.Range("Q367").Value = 12

# Row 368
$ws.Range("D368").Value = 44169
$ws.Range("D368").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H368").Value = "Inferno"
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 150
$ws.Range("K368").Value = 28000
$ws.Range("L368").Value = 28000
$ws.Range("M368").Value = 28000
$ws.Range("N368").Value = "`$/caja 12 kilos"
$ws.Range("O368").Value = "Región de Arica y Parinacota"
$ws.Range("P368").Value = 2333
$ws.Range("Q368").Value = 12

# Row 369
$ws.Range("D369").Value = 44351
$ws.Range("D369").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H369").Value = "Inferno"
$ws.Range("I369").Value = "Primera"
$ws.Range("J369").Value = 80
$ws.Range("K369").Value = 31000
$ws.Range("L369").Value = 31000
$ws.Range("M369").Value = 31000
$ws.Range("N369").Value = "`$/caja 12 kilos"
$ws.Range("O369").Value = "Región de Arica y Parinacota"
$ws.Range("P369").Value = 2583
$ws.Range("Q369").Value = 12

# Row 370
$ws.Range("D370").Value = 44795
$ws.Range("D370").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H370").Value = "Inferno"
$ws.Range("I370").Value = "Primera"
$ws.Range("J370").Value = 70
$ws.Range("K370").Value = 19000
$ws.Range("L370").Value = 19000
$ws.Range("M370").Value = 19000
$ws.Range("N370").Value = "`$/caja 10 kilos"
$ws.Range("O370").Value = "Región de Arica y Parinacota"
$ws.Range("P370").Value = 1900
$ws.Range("Q370").Value = 10

# Row 371
$ws.Range("D371").Value = 44442
$ws.Range("D371").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H371").Value = "Inferno"
$ws.Range("I371").Value = "Primera"
$ws.Range("J371").Value = 120
$ws.Range("K371").Value = 45000
$ws.Range("L371").Value = 45000
$ws.Range("M371").Value = 45000
$ws.Range("N371").Value = "`$/caja 12 kilos"
$ws.Range("O371").Value = "Región de Arica y Parinacota"
$ws.Range("P371").Value = 3750
$ws.Range("Q371").Value = 12

# Row 372
$ws.Range("D372").Value = 44595
$ws.Range("D372").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H372").Value = "Inferno"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 50
$ws.Range("K372").Value = 17000
$ws.Range("L372").Value = 17000
$ws.Range("M372").Value = 17000
$ws.Range("N372").Value = "`$/caja 12 kilos"
$ws.Range("O372").Value = "Región de Arica y Parinacota"
$ws.Range("P372").Value = 1417
$ws.Range("Q372").Value = 12
